$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Update "Your algorithm best solution obj" values (column I) for a subset
# of instances - this drives the dependent Gap% formulas in column J and
# the summary AVERAGE in J13 to recalculate automatically.
$ws.Range("I5").Value = 39.9
$ws.Range("I6").Value = 39.14
$ws.Range("I7").Value = 9.15
$ws.Range("I10").Value = 10.957000000000001

# Move the active selection to I11
$ws.Range("I11").Select()
